$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 8-10 (data now only spans rows 2-7)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Hgf"
$ws.Range("C2").Value = "Sdc2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.583520999999999
$ws.Range("H2").Value = 25.750563
$ws.Range("I2").Value = 0.8910607110509009
$ws.Range("J2").Value = 0.8910607110509009
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.116695
$ws.Range("N2").Value = 3.350085
$ws.Range("O2").Value = 0.008174214292497491
$ws.Range("P2").Value = 0.008174214292497492
$ws.Range("Q2").Value = 9.585174983094999
$ws.Range("R2").Value = 86.266574847855
$ws.Range("S2").Value = 0.00728372119975525
$ws.Range("T2").Value = 0.007283721199755252

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Hgf"
$ws.Range("C3").Value = "Sdc2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.583520999999999
$ws.Range("H3").Value = 25.750563
$ws.Range("I3").Value = 0.8910607110509009
$ws.Range("J3").Value = 0.8910607110509009
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 111.9320066666667
$ws.Range("N3").Value = 335.79602
$ws.Range("O3").Value = 0.8193429796700005
$ws.Range("P3").Value = 0.8193429796700005
$ws.Range("Q3").Value = 960.7707297954732
$ws.Range("R3").Value = 8646.936568159259
$ws.Range("S3").Value = 0.7300843380593145
$ws.Range("T3").Value = 0.7300843380593145

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Hgf"
$ws.Range("C4").Value = "Sdc2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.583520999999999
$ws.Range("H4").Value = 25.750563
$ws.Range("I4").Value = 0.8910607110509009
$ws.Range("J4").Value = 0.8910607110509009
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.563205
$ws.Range("N4").Value = 70.689615
$ws.Range("O4").Value = 0.172482806037502
$ws.Range("P4").Value = 0.1724828060375021
$ws.Range("Q4").Value = 202.255264944805
$ws.Range("R4").Value = 1820.297384503245
$ws.Range("S4").Value = 0.1536926517918312
$ws.Range("T4").Value = 0.1536926517918312

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Hgf"
$ws.Range("C5").Value = "Sdc2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.049404
$ws.Range("H5").Value = 3.148212
$ws.Range("I5").Value = 0.1089392889490991
$ws.Range("J5").Value = 0.1089392889490991
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.116695
$ws.Range("N5").Value = 3.350085
$ws.Range("O5").Value = 0.008174214292497491
$ws.Range("P5").Value = 0.008174214292497492
$ws.Range("Q5").Value = 1.17186419978
$ws.Range("R5").Value = 10.54677779802
$ws.Range("S5").Value = 0.0008904930927422395
$ws.Range("T5").Value = 0.0008904930927422397

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Hgf"
$ws.Range("C6").Value = "Sdc2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.049404
$ws.Range("H6").Value = 3.148212
$ws.Range("I6").Value = 0.1089392889490991
$ws.Range("J6").Value = 0.1089392889490991
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 111.9320066666667
$ws.Range("N6").Value = 335.79602
$ws.Range("O6").Value = 0.8193429796700005
$ws.Range("P6").Value = 0.8193429796700005
$ws.Range("Q6").Value = 117.4618955240267
$ws.Range("R6").Value = 1057.15705971624
$ws.Range("S6").Value = 0.08925864161068597
$ws.Range("T6").Value = 0.08925864161068597

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Hgf"
$ws.Range("C7").Value = "Sdc2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.049404
$ws.Range("H7").Value = 3.148212
$ws.Range("I7").Value = 0.1089392889490991
$ws.Range("J7").Value = 0.1089392889490991
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.563205
$ws.Range("N7").Value = 70.689615
$ws.Range("O7").Value = 0.172482806037502
$ws.Range("P7").Value = 0.1724828060375021
$ws.Range("Q7").Value = 24.72732157982
$ws.Range("R7").Value = 222.54589421838
$ws.Range("S7").Value = 0.01879015424567084
$ws.Range("T7").Value = 0.01879015424567084

